# Add new test cases for "multiple routing" to the snc-connector test data
# sheet (getConceptModelDataByCondition).
#
# The existing block of rows describing the "differentTablesInTwoDb"
# scenarios (rows 81-85) gets re-arranged so that the "order" / "fields" /
# "condition" example values are shifted down by one row, row 83 gets a
# brand-new "order" example, and three new rows (86, 87, 88) are appended
# with new scenarios covering multiple routing.
#
# NOTE: cell writes below are intentionally ordered so that brand-new text
# values are appended to the shared-strings table in the same order as the
# original authoring session (this keeps the generated workbook's
# xl/sharedStrings.xml index layout identical to the target).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copies the border/format of an existing "description" cell (column B) onto
# a new cell, so the new cell keeps style index 3 (bordered) instead of the
# bare default style that brand-new cells otherwise get.
function Set-DescCell($row) {
    $ws.Cells.Item(81, 2).Copy() | Out-Null
    $ws.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null
    $ws.Application.CutCopyMode() = 0
}

# ---------------------------------------------------------------------------
# Row 86 (new): re-create the "order"/pageIndex/pageSize example that used to
# live on row 85, under a new test-id.
# ---------------------------------------------------------------------------
$ws.Cells.Item(86, 1).Value() = "snc-connector-differentTablesInTwoDb-6"
Set-DescCell 86
$ws.Cells.Item(86, 2).Value() = "good request, data retrieved (no schema check)"
$ws.Cells.Item(86, 6).Value() = "Work_Center3"
$ws.Cells.Item(86, 7).Value() = "location"
$ws.Cells.Item(86, 8).Value() = 1
$ws.Cells.Item(86, 9).Value() = 2
$ws.Cells.Item(86, 12).Value() = 0

# ---------------------------------------------------------------------------
# Row 87 (new): multiple-routing bad-request case - the request condition
# includes the rule key.
# ---------------------------------------------------------------------------
$ws.Cells.Item(87, 3).Value() = "plant_owner=3"

# ---------------------------------------------------------------------------
# Row 83: drop the "condition" value that used to live here, add a new
# "order" value instead.
# ---------------------------------------------------------------------------
$ws.Cells.Item(83, 3).Clear() | Out-Null
$ws.Cells.Item(83, 7).Value() = "location,-description"

Set-DescCell 87
$ws.Cells.Item(87, 2).Value() = "bad request (condition includes rule key)"

# ---------------------------------------------------------------------------
# Row 88 (new): multiple-routing bad-request case - the generated SQL
# execution itself fails.
# ---------------------------------------------------------------------------
$ws.Cells.Item(88, 13).Value() = "sql execution failed:Error while executing SQL"

$ws.Cells.Item(87, 1).Value() = "snc-connector-differentTablesInTwoDb-7"
$ws.Cells.Item(88, 1).Value() = "snc-connector-differentTablesInTwoDb-8"

$ws.Cells.Item(87, 13).Value() = "do not support ruleKey value"

$ws.Cells.Item(87, 6).Value() = "Work_Center3"
$ws.Cells.Item(87, 12).Value() = 106601

Set-DescCell 88
$ws.Cells.Item(88, 2).Value() = "bad request (condition includes rule key)"
$ws.Cells.Item(88, 6).Value() = "Work_Center3"
$ws.Cells.Item(88, 8).Value() = 1
$ws.Cells.Item(88, 9).Value() = 2
$ws.Cells.Item(88, 12).Value() = 106103

# ---------------------------------------------------------------------------
# Row 84: drop the "fields" value, add the "condition" value that used to be
# on row 83 (id='A5E03262697').
# ---------------------------------------------------------------------------
$ws.Cells.Item(84, 5).Clear() | Out-Null
$ws.Cells.Item(84, 3).Value() = "id='A5E03262697'"

# ---------------------------------------------------------------------------
# Row 85: drop the "order"/pageIndex/pageSize values, add the "fields" value
# that used to be on row 84 (description).
# ---------------------------------------------------------------------------
$ws.Cells.Item(85, 7).Clear() | Out-Null
$ws.Cells.Item(85, 8).Clear() | Out-Null
$ws.Cells.Item(85, 9).Clear() | Out-Null
$ws.Cells.Item(85, 5).Value() = "description"

# ---------------------------------------------------------------------------
# View tweaks: widen column M (rspMessage) to match column C since it now
# holds message text, and move the selection down to the new last row.
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth() = $ws.Columns.Item(3).ColumnWidth()

$ws.Range("B89").Select() | Out-Null
